$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete trailing rows (old rows 23 and 24)
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Delete()

# Update header row: new Jan_2026 column inserted, Oct_2025 column dropped
$ws.Cells.Item(1,4).Value = 'Jan_2026'
$ws.Cells.Item(1,5).Value = 'Dec_2025'
$ws.Cells.Item(1,6).Value = 'Nov_2025'

# Rewrite holdings rows 2-22 with refreshed quant engine data
$ws.Cells.Item(2,1).Value = 'INE040A01034'
$ws.Cells.Item(2,2).Value = 'HDFC Bank Limited'
$ws.Cells.Item(2,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(2,4).Value = 9.269297
$ws.Cells.Item(2,5).Value = 5.688921
$ws.Cells.Item(2,6).Value = 5.548269
$ws.Cells.Item(2,7).Value = 3.580376
$ws.Cells.Item(2,8).Value = 3.721028

$ws.Cells.Item(3,1).Value = 'INE090A01021'
$ws.Cells.Item(3,2).Value = 'ICICI Bank Limited'
$ws.Cells.Item(3,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(3,4).Value = 8.668876
$ws.Cells.Item(3,5).Value = 4.458872
$ws.Cells.Item(3,6).Value = 6.814643
$ws.Cells.Item(3,7).Value = 4.210003999999999
$ws.Cells.Item(3,8).Value = 1.854232999999999

$ws.Cells.Item(4,1).Value = 'INE364U01010'
$ws.Cells.Item(4,2).Value = 'Adani Green Energy Limited'
$ws.Cells.Item(4,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(4,4).Value = 7.780966
$ws.Cells.Item(4,5).Value = 9.24471
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = -1.463743999999999
$ws.Cells.Item(4,8).Value = 7.780966

$ws.Cells.Item(5,1).Value = 'INE397D01024'
$ws.Cells.Item(5,2).Value = 'Bharti Airtel Limited'
$ws.Cells.Item(5,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(5,4).Value = 6.797444
$ws.Cells.Item(5,5).Value = 5.218497
$ws.Cells.Item(5,6).Value = 4.997126
$ws.Cells.Item(5,7).Value = 1.578946999999999
$ws.Cells.Item(5,8).Value = 1.800318

$ws.Cells.Item(6,1).Value = 'INE326A01037'
$ws.Cells.Item(6,2).Value = 'Lupin Limited'
$ws.Cells.Item(6,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(6,4).Value = 5.002997
$ws.Cells.Item(6,5).Value = 4.890458
$ws.Cells.Item(6,6).Value = 4.631196
$ws.Cells.Item(6,7).Value = 0.1125389999999999
$ws.Cells.Item(6,8).Value = 0.3718009999999996

$ws.Cells.Item(7,1).Value = 'INE296A01032'
$ws.Cells.Item(7,2).Value = 'Bajaj Finance Limited'
$ws.Cells.Item(7,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(7,4).Value = 4.973993
$ws.Cells.Item(7,5).Value = 5.265807
$ws.Cells.Item(7,6).Value = 5.311591
$ws.Cells.Item(7,7).Value = -0.2918139999999996
$ws.Cells.Item(7,8).Value = -0.3375979999999998

$ws.Cells.Item(8,1).Value = 'INE002A01018'
$ws.Cells.Item(8,2).Value = 'Reliance Industries Limited'
$ws.Cells.Item(8,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(8,4).Value = 4.644474
$ws.Cells.Item(8,5).Value = 5.214247
$ws.Cells.Item(8,6).Value = 5.978789
$ws.Cells.Item(8,7).Value = -0.5697730000000005
$ws.Cells.Item(8,8).Value = -1.334315

$ws.Cells.Item(9,1).Value = 'INE061F01013'
$ws.Cells.Item(9,2).Value = 'Fortis Healthcare Ltd'
$ws.Cells.Item(9,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(9,4).Value = 4.392929
$ws.Cells.Item(9,5).Value = 4.549517
$ws.Cells.Item(9,6).Value = 4.538125
$ws.Cells.Item(9,7).Value = -0.1565880000000002
$ws.Cells.Item(9,8).Value = -0.1451960000000003

$ws.Cells.Item(10,1).Value = 'INE020B01018'
$ws.Cells.Item(10,2).Value = 'Rural Electrification Corporation Ltd'
$ws.Cells.Item(10,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(10,4).Value = 4.362757
$ws.Cells.Item(10,5).Value = 4.264898
$ws.Cells.Item(10,6).Value = 4.083341
$ws.Cells.Item(10,7).Value = 0.09785900000000058
$ws.Cells.Item(10,8).Value = 0.2794160000000003

$ws.Cells.Item(11,1).Value = 'INE303R01014'
$ws.Cells.Item(11,2).Value = 'Kalyan Jewellers India Limited'
$ws.Cells.Item(11,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(11,4).Value = 4.344148
$ws.Cells.Item(11,5).Value = 5.815869
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = -1.471721000000001
$ws.Cells.Item(11,8).Value = 4.344148

$ws.Cells.Item(12,1).Value = 'INE758T01015'
$ws.Cells.Item(12,2).Value = 'Eternal Limited'
$ws.Cells.Item(12,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(12,4).Value = 4.258939
$ws.Cells.Item(12,5).Value = 4.317692
$ws.Cells.Item(12,6).Value = 4.470905
$ws.Cells.Item(12,7).Value = -0.05875300000000028
$ws.Cells.Item(12,8).Value = -0.2119660000000003

$ws.Cells.Item(13,1).Value = 'INE584A01023'
$ws.Cells.Item(13,2).Value = 'NMDC Ltd'
$ws.Cells.Item(13,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(13,4).Value = 1.954854
$ws.Cells.Item(13,5).Value = 1.99717
$ws.Cells.Item(13,6).Value = 4.541295
$ws.Cells.Item(13,7).Value = -0.0423159999999998
$ws.Cells.Item(13,8).Value = -2.586441

$ws.Cells.Item(14,1).Value = 'INE271C01023'
$ws.Cells.Item(14,2).Value = 'DLF Limited'
$ws.Cells.Item(14,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 2.575401
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = -2.575401

$ws.Cells.Item(15,1).Value = 'INE237A01028'
$ws.Cells.Item(15,2).Value = 'Kotak Mahindra Bank Limited'
$ws.Cells.Item(15,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 6.264325
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = -6.264325
$ws.Cells.Item(15,8).Value = 0

$ws.Cells.Item(16,1).Value = 'INE115A01026'
$ws.Cells.Item(16,2).Value = 'LIC Housing Finance Ltd'
$ws.Cells.Item(16,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 2.655749
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = -2.655749

$ws.Cells.Item(17,1).Value = 'INE0BS701011'
$ws.Cells.Item(17,2).Value = 'Premier Energies Limited'
$ws.Cells.Item(17,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 5.023172
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = -5.023172

$ws.Cells.Item(18,1).Value = 'INE437A01024'
$ws.Cells.Item(18,2).Value = 'Apollo Hospitals Enterprise Ltd'
$ws.Cells.Item(18,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 1.50322
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = -1.50322
$ws.Cells.Item(18,8).Value = 0

$ws.Cells.Item(19,1).Value = 'INE467B01029'
$ws.Cells.Item(19,2).Value = 'Tata Consultancy Services Limited'
$ws.Cells.Item(19,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 4.92361
$ws.Cells.Item(19,6).Value = 4.622506
$ws.Cells.Item(19,7).Value = -4.92361
$ws.Cells.Item(19,8).Value = -4.622506

$ws.Cells.Item(20,1).Value = 'INE484J01027'
$ws.Cells.Item(20,2).Value = 'Godrej Properties Limited'
$ws.Cells.Item(20,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 2.465903
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = -2.465903

$ws.Cells.Item(21,1).Value = 'INE647A01010'
$ws.Cells.Item(21,2).Value = 'SRF Limited'
$ws.Cells.Item(21,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 1.264205
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = -1.264205

$ws.Cells.Item(22,1).Value = 'INE917I01010'
$ws.Cells.Item(22,2).Value = 'Bajaj Auto Limited'
$ws.Cells.Item(22,3).Value = 'quant Equity Savings Fund'
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 3.057868
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = -3.057868
$ws.Cells.Item(22,8).Value = 0

